$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The uploaded data set dropped one category row (the old B27 = 486542
# entry). Every row below it shifts up by one (columns B:D only -- the
# running index in column A stays untouched), and the now-duplicate last
# row (70) disappears entirely.
for ($r = 27; $r -le 69; $r++) {
    $srcRow = $r + 1
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($srcRow, 2).Value2
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($srcRow, 3).Value2
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($srcRow, 4).Value2
}

# Drop the now-redundant trailing row; this also shrinks the sheet
# dimension from A1:D70 to A1:D69 automatically.
$ws.Rows.Item(70).Delete()

# Match the author's saved cursor position.
$ws.Range("D2").Select() | Out-Null
